$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="D2"; Value="63.358.91"; AsText=$false},
    @{Cell="E2"; Value="  -1.13%  "; AsText=$false},
    @{Cell="D3"; Value="2.683.93"; AsText=$false},
    @{Cell="E3"; Value="  -2.78%  "; AsText=$false},
    @{Cell="E4"; Value="  +0.00%  "; AsText=$false},
    @{Cell="D5"; Value="553.82"; AsText=$true},
    @{Cell="E5"; Value="  -3.73%  "; AsText=$false},
    @{Cell="D6"; Value="158.37"; AsText=$true},
    @{Cell="E6"; Value="  -0.36%  "; AsText=$false},
    @{Cell="D7"; Value="0.999"; AsText=$true},
    @{Cell="E7"; Value="  +0.05%  "; AsText=$false},
    @{Cell="D8"; Value="0.585"; AsText=$true},
    @{Cell="E8"; Value="  -2.75%  "; AsText=$false},
    @{Cell="E9"; Value="  -3.87%  "; AsText=$false},
    @{Cell="E10"; Value="  -0.66%  "; AsText=$false},
    @{Cell="D12"; Value="5.34"; AsText=$true},
    @{Cell="E12"; Value="  -9.12%  "; AsText=$false},
    @{Cell="D13"; Value="3.155.72"; AsText=$false},
    @{Cell="E13"; Value="  -2.90%  "; AsText=$false},
    @{Cell="D14"; Value="26.37"; AsText=$true},
    @{Cell="E14"; Value="  -2.06%  "; AsText=$false},
    @{Cell="D15"; Value="63.220.96"; AsText=$false},
    @{Cell="E15"; Value="  -0.81%  "; AsText=$false},
    @{Cell="E16"; Value="  -4.07%  "; AsText=$false},
    @{Cell="D17"; Value="2.683.05"; AsText=$false},
    @{Cell="E17"; Value="  -3.13%  "; AsText=$false},
    @{Cell="D18"; Value="11.99"; AsText=$true},
    @{Cell="E18"; Value="  -1.44%  "; AsText=$false},
    @{Cell="D19"; Value="4.56"; AsText=$true},
    @{Cell="E19"; Value="  -5.47%  "; AsText=$false},
    @{Cell="D20"; Value="344.57"; AsText=$true},
    @{Cell="E20"; Value="  -4.30%  "; AsText=$false},
    @{Cell="E21"; Value="  -4.79%  "; AsText=$false},
    @{Cell="D22"; Value="0.996"; AsText=$true},
    @{Cell="E22"; Value="  -0.27%  "; AsText=$false},
    @{Cell="E23"; Value="  -3.87%  "; AsText=$false},
    @{Cell="D24"; Value="63.75"; AsText=$true},
    @{Cell="E24"; Value="  -1.93%  "; AsText=$false},
    @{Cell="D25"; Value="0.169"; AsText=$true},
    @{Cell="E25"; Value="  -1.04%  "; AsText=$false},
    @{Cell="D26"; Value="1.00"; AsText=$true},
    @{Cell="E26"; Value="  +0.31%  "; AsText=$false},
    @{Cell="D27"; Value="8.17"; AsText=$true},
    @{Cell="E27"; Value="  -4.18%  "; AsText=$false},
    @{Cell="D28"; Value="0.0₃0854"; AsText=$false},
    @{Cell="E28"; Value="  -5.54%  "; AsText=$false},
    @{Cell="D29"; Value="1.94"; AsText=$true},
    @{Cell="E29"; Value="  -1.13%  "; AsText=$false},
    @{Cell="E30"; Value="  +0.43%  "; AsText=$false},
    @{Cell="D31"; Value="7.02"; AsText=$true},
    @{Cell="E31"; Value="  -4.72%  "; AsText=$false},
    @{Cell="D32"; Value="165.91"; AsText=$true},
    @{Cell="E32"; Value="  -1.76%  "; AsText=$false},
    @{Cell="E33"; Value="  +0.02%  "; AsText=$false},
    @{Cell="D34"; Value="4.82"; AsText=$true},
    @{Cell="E34"; Value="  -2.52%  "; AsText=$false},
    @{Cell="D35"; Value="19.51"; AsText=$true},
    @{Cell="E35"; Value="  -3.28%  "; AsText=$false},
    @{Cell="E36"; Value="  -4.76%  "; AsText=$false},
    @{Cell="E37"; Value="  -1.71%  "; AsText=$false},
    @{Cell="D38"; Value="341.72"; AsText=$true},
    @{Cell="E38"; Value="  -2.20%  "; AsText=$false},
    @{Cell="D39"; Value="0.944"; AsText=$true},
    @{Cell="E39"; Value="  -5.95%  "; AsText=$false},
    @{Cell="D40"; Value="6.10"; AsText=$true},
    @{Cell="E40"; Value="  -3.28%  "; AsText=$false},
    @{Cell="D41"; Value="38.04"; AsText=$true},
    @{Cell="E41"; Value="  -2.67%  "; AsText=$false},
    @{Cell="D42"; Value="3.92"; AsText=$true},
    @{Cell="E42"; Value="  -6.37%  "; AsText=$false},
    @{Cell="B43"; Value="InjectiveProtocol"; AsText=$false},
    @{Cell="C43"; Value="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; AsText=$false},
    @{Cell="D43"; Value="20.75"; AsText=$true},
    @{Cell="E43"; Value="  -6.20%  "; AsText=$false},
    @{Cell="B44"; Value="EnergySwap"; AsText=$false},
    @{Cell="C44"; Value="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; AsText=$false},
    @{Cell="D44"; Value="20.25"; AsText=$true},
    @{Cell="E44"; Value="  -5.95%  "; AsText=$false},
    @{Cell="D45"; Value="0.618"; AsText=$true},
    @{Cell="E45"; Value="  -1.58%  "; AsText=$false},
    @{Cell="D46"; Value="0.0562"; AsText=$true},
    @{Cell="E46"; Value="  -4.54%  "; AsText=$false},
    @{Cell="D47"; Value="0.998"; AsText=$true},
    @{Cell="E47"; Value="  +0.01%  "; AsText=$false},
    @{Cell="D48"; Value="11.06"; AsText=$true},
    @{Cell="E48"; Value="  +0.17%  "; AsText=$false},
    @{Cell="D49"; Value="129.63"; AsText=$true},
    @{Cell="E49"; Value="  -6.01%  "; AsText=$false},
    @{Cell="D50"; Value="0.0971"; AsText=$true},
    @{Cell="E50"; Value="  -3.75%  "; AsText=$false},
    @{Cell="E51"; Value="  -4.74%  "; AsText=$false}
)

foreach ($u in $updates) {
    if ($u.AsText) {
        $ws.Range($u.Cell).NumberFormat = "@"
    }
    $ws.Range($u.Cell).Value = $u.Value
}
